$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 2 (Saurav) data - bug fix: hours worked and active status
$ws.Range("G2").Value = 19.0
$ws.Range("I2").Value = $false

# Add new manager rows (fixing bug with item/manager creation)
$ws.Range("A4").Value = 11.0
$ws.Range("B4").Value = "testmanager@gmail.com"
$ws.Range("C4").Value = "Test"
$ws.Range("D4").Value = "Manager"
$ws.Range("E4").Value = 6.0
$ws.Range("F4").Value = "Manager"
$ws.Range("G4").Value = 0.0
$ws.Range("H4").Value = 20.0
$ws.Range("I4").Value = $true

$ws.Range("A5").Value = 17.0
$ws.Range("B5").Value = "manager@gmail.com"
$ws.Range("C5").Value = "Manager"
$ws.Range("D5").Value = "Test"
$ws.Range("E5").Value = 11.0
$ws.Range("F5").Value = "Manager"
$ws.Range("G5").Value = 0.0
$ws.Range("H5").Value = 20.0
$ws.Range("I5").Value = $true
